# Fix issue mapping Medication ab77bff992840669d45583ae812eee5175aff7fe
$wb = $excel.ActiveWorkbook

$wsMapping = $wb.Worksheets.Item("Mapping Table 0")

# Duplicate row 5's formatting down into a new row 6 so the new mapping row
# picks up the same styling (borders/fill/alignment) as the existing rows.
$wsMapping.Range("A5:E5").Copy()
$wsMapping.Range("A6:E6").PasteSpecial(-4122)  # xlPasteFormats

# The duplicated "Forme" Source entries on rows 4 and 5 no longer apply
# directly - blank out their Source column.
$wsMapping.Range("A4").Value = ""
$wsMapping.Range("A5").Value = ""

# Add the new "related-to" relationship row, still targeting
# Medication.form.coding.code, with blank Source/Display columns.
$wsMapping.Range("C6").Value = "related-to"
$wsMapping.Range("D6").Value = "Medication.form.coding.code"

# Update the metadata "Date" value on the Metadata sheet to reflect the edit time.
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2024-11-17T10:38:58+00:00"
